# Fruta / hortaliza, semanal
# Insert 3 new weekly observations for Palta (Hass) at rows 257-259,
# pushing the previously existing rows 257-290 down to 260-293.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 257 (shifts old rows 257:290 -> 260:293)
$ws.Range("A257:A259").EntireRow.Insert()

# Constant columns shared by every "Palta" row on this sheet
$mercadoId  = 7
$mercado    = "Terminal Hortofrutícola Agro Chillán"
$region     = "Ñuble"
$codreg     = 16
$tipo       = "Fruta"
$productoId = 100106
$producto   = "Oleaginosos"
$categoriaId = 100106002
$categoria  = "Palta"
$variedad   = "Hass"

function Set-PaltaRow {
    param(
        [int]$row,
        [double]$fecha,
        [string]$calidad,
        [double]$volumen,
        [double]$precioMin,
        [double]$precioMax,
        [double]$precioProm,
        [string]$unidad,
        [string]$origen,
        [double]$precioKg,
        [double]$kgUnidad
    )

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

# New row 257: 1a nueva(o), Provincia de Quillota
Set-PaltaRow 257 44474 "1a nueva(o)" 120 2900 3000 2950 "$/kilo (en caja de 15 kilos)" "Provincia de Quillota" 2950 1

# New row 258: Primera, Perú
Set-PaltaRow 258 44474 "Primera" 400 25000 26000 25500 "$/bandeja 10 kilos" "Perú" 2550 10

# New row 259: Segunda, Perú
Set-PaltaRow 259 44474 "Segunda" 240 23000 24000 23500 "$/bandeja 10 kilos" "Perú" 2350 10
